# Auto-generated edit script: updates crypto price/volume figures
# and fixes the EnergySwap/Decentraland row ordering, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.978.60' }
    @{ Cell = 'D3'; Value = '1.880.61' }
    @{ Cell = 'E4'; Value = '  +0.11%  ' }
    @{ Cell = 'D5'; Value = '''242.64' }
    @{ Cell = 'E5'; Value = '  -3.96%  ' }
    @{ Cell = 'D6'; Value = '''0.9999' }
    @{ Cell = 'E6'; Value = '  +0.05%  ' }
    @{ Cell = 'D7'; Value = '''0.4917' }
    @{ Cell = 'E7'; Value = '  -3.41%  ' }
    @{ Cell = 'D8'; Value = '''0.2940' }
    @{ Cell = 'E8'; Value = '  -2.86%  ' }
    @{ Cell = 'D9'; Value = '''0.06615' }
    @{ Cell = 'E9'; Value = '  -3.11%  ' }
    @{ Cell = 'D10'; Value = '1.885.72' }
    @{ Cell = 'E10'; Value = '  -1.17%  ' }
    @{ Cell = 'D11'; Value = '''16.67' }
    @{ Cell = 'E11'; Value = '  -3.66%  ' }
    @{ Cell = 'E12'; Value = '  -2.08%  ' }
    @{ Cell = 'D13'; Value = '''0.6660' }
    @{ Cell = 'E13'; Value = '  -3.77%  ' }
    @{ Cell = 'D14'; Value = '''86.30' }
    @{ Cell = 'E14'; Value = '  -0.75%  ' }
    @{ Cell = 'D15'; Value = '''4.872' }
    @{ Cell = 'E15'; Value = '  -0.96%  ' }
    @{ Cell = 'D16'; Value = '29.975.64' }
    @{ Cell = 'E16'; Value = '  -0.47%  ' }
    @{ Cell = 'D17'; Value = '''0.000007808' }
    @{ Cell = 'E17'; Value = '  -6.12%  ' }
    @{ Cell = 'E18'; Value = '  +0.14%  ' }
    @{ Cell = 'D19'; Value = '''12.78' }
    @{ Cell = 'E19'; Value = '  -2.15%  ' }
    @{ Cell = 'D20'; Value = '2.124.21' }
    @{ Cell = 'E20'; Value = '  -1.29%  ' }
    @{ Cell = 'D21'; Value = '''1.001' }
    @{ Cell = 'E21'; Value = '  +0.20%  ' }
    @{ Cell = 'D22'; Value = '''4.776' }
    @{ Cell = 'E22'; Value = '  -1.01%  ' }
    @{ Cell = 'D23'; Value = '''5.845' }
    @{ Cell = 'E23'; Value = '  +1.84%  ' }
    @{ Cell = 'D24'; Value = '''9.080' }
    @{ Cell = 'E24'; Value = '  -2.31%  ' }
    @{ Cell = 'D25'; Value = '''150.52' }
    @{ Cell = 'E25'; Value = '  +2.12%  ' }
    @{ Cell = 'D26'; Value = '''140.71' }
    @{ Cell = 'E26'; Value = '  +4.57%  ' }
    @{ Cell = 'D27'; Value = '''16.94' }
    @{ Cell = 'E27'; Value = '  -1.18%  ' }
    @{ Cell = 'D28'; Value = '''1.898' }
    @{ Cell = 'E28'; Value = '  -5.39%  ' }
    @{ Cell = 'D29'; Value = '''1.390' }
    @{ Cell = 'E29'; Value = '  -0.59%  ' }
    @{ Cell = 'D30'; Value = '''4.186' }
    @{ Cell = 'E30'; Value = '  -2.32%  ' }
    @{ Cell = 'D31'; Value = '''0.08747' }
    @{ Cell = 'E31'; Value = '  -1.33%  ' }
    @{ Cell = 'E32'; Value = '  -0.62%  ' }
    @{ Cell = 'D33'; Value = '''0.05014' }
    @{ Cell = 'E33'; Value = '  -0.90%  ' }
    @{ Cell = 'D34'; Value = '''0.7158' }
    @{ Cell = 'E34'; Value = '  -1.06%  ' }
    @{ Cell = 'D35'; Value = '''1.109' }
    @{ Cell = 'E35'; Value = '  -3.16%  ' }
    @{ Cell = 'D36'; Value = '''2.669' }
    @{ Cell = 'E36'; Value = '  -0.76%  ' }
    @{ Cell = 'D37'; Value = '''0.01787' }
    @{ Cell = 'E37'; Value = '  +5.50%  ' }
    @{ Cell = 'D38'; Value = '''2.693' }
    @{ Cell = 'E38'; Value = '  -4.66%  ' }
    @{ Cell = 'D39'; Value = '''2.157' }
    @{ Cell = 'E39'; Value = '  -5.22%  ' }
    @{ Cell = 'D40'; Value = '''0.9397' }
    @{ Cell = 'E40'; Value = '  -2.33%  ' }
    @{ Cell = 'D41'; Value = '''0.9995' }
    @{ Cell = 'E41'; Value = '  +0.05%  ' }
    @{ Cell = 'D42'; Value = '''103.68' }
    @{ Cell = 'E42'; Value = '  -1.13%  ' }
    @{ Cell = 'D43'; Value = '''0.4225' }
    @{ Cell = 'E43'; Value = '  -2.05%  ' }
    @{ Cell = 'D44'; Value = '''5.740' }
    @{ Cell = 'E44'; Value = '  -6.26%  ' }
    @{ Cell = 'D45'; Value = '''7.318' }
    @{ Cell = 'E45'; Value = '  -4.64%  ' }
    @{ Cell = 'D46'; Value = '''0.1268' }
    @{ Cell = 'E46'; Value = '  -0.95%  ' }
    @{ Cell = 'D47'; Value = '''0.05698' }
    @{ Cell = 'E47'; Value = '  -0.88%  ' }
    @{ Cell = 'D48'; Value = '''32.66' }
    @{ Cell = 'E48'; Value = '  -1.74%  ' }
    @{ Cell = 'B49'; Value = 'Decentraland' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D49'; Value = '''0.3752' }
    @{ Cell = 'E49'; Value = '  -1.97%  ' }
    @{ Cell = 'B50'; Value = 'EnergySwap' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D50'; Value = '''8.215' }
    @{ Cell = 'E50'; Value = '  -2.27%  ' }
    @{ Cell = 'D51'; Value = '''55.95' }
    @{ Cell = 'E51'; Value = '  -1.65%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
